$wb = $excel.ActiveWorkbook

$index = $wb.Worksheets.Item("Index")

# Add the two new rows with the same style as existing hyperlink rows
$index.Hyperlinks.Add($index.Range("A3"), "", "child1!A1", "", "child1")
$index.Hyperlinks.Add($index.Range("A4"), "", "child2!A1", "", "child2")

$index.Range("A3:A4").Style = $index.Range("A2").Style

$index.Range("A4").Select()

# Make Index the active sheet / tab
$index.Activate()

$wb.Save()
